# Add a new wave of survey data (28. 9. 2021) as an additional column
# to both worksheets: "data" (sheet1, new column AI) and "pocetR" (sheet2, new column AH).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "data" -> new column AI, header date "28. 9. 2021"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("data")

# Copy formatting of the previous header cell (AH1) into the new header cell (AI1)
$ws1.Range("AH1").Copy()
$ws1.Range("AI1").PasteSpecial(-4122)
$ws1.Range("AI1").Value = "28. 9. 2021"

# Percentage values for rows 2..76 for the new wave
$ws1Values = @{
    2  = 0.29
    3  = 0.43
    4  = 0.28
    5  = 0.38
    6  = 0.44
    7  = 0.18
    8  = 0.3
    9  = 0.46
    10 = 0.24
    11 = 0.23
    12 = 0.4
    13 = 0.37
    14 = 0.31
    15 = 0.4
    16 = 0.29
    17 = 0.27
    18 = 0.4
    19 = 0.33
    20 = 0.29
    21 = 0.45
    22 = 0.26
    23 = 0.36
    24 = 0.41
    25 = 0.23
    26 = 0.23
    27 = 0.45
    28 = 0.32
    29 = 0.28
    30 = 0.43
    31 = 0.29
    32 = 0.31
    33 = 0.44
    34 = 0.25
    35 = 0.32
    36 = 0.4
    37 = 0.28
    38 = 0.27
    39 = 0.45
    40 = 0.28
    41 = 0.11
    42 = 0.44
    43 = 0.45
    44 = 0.35
    45 = 0.33
    46 = 0.32
    47 = 0.31
    48 = 0.46
    49 = 0.23
    50 = 0.3
    51 = 0.5600000000000001
    52 = 0.14
    53 = 0.37
    54 = 0.44
    55 = 0.19
    56 = 0.32
    57 = 0.47
    58 = 0.21
    59 = 0.44
    60 = 0.41
    61 = 0.15
    62 = 0.39
    63 = 0.41
    64 = 0.2
    65 = 0.27
    66 = 0.41
    67 = 0.32
    68 = 0.32
    69 = 0.47
    70 = 0.21
    71 = 0.2
    72 = 0.51
    73 = 0.29
    74 = 0.19
    75 = 0.39
    76 = 0.42
}

for ($row = 2; $row -le 76; $row++) {
    $ws1.Range("AI$row").Value = $ws1Values[$row]
}

# Update the footnote text in row 77 with the new update date
$ws1.Range("A77").Value = "Život během pandemie, Obavy z epidemie, % respondentů celkově a ve skupinách, aktualizace 6. 10. 2021"

# ---------------------------------------------------------------------------
# Sheet 2: "pocetR" -> new column AH, header date "28. 9. 2021"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("pocetR")

# Copy formatting of the previous header cell (AG1) into the new header cell (AH1)
$ws2.Range("AG1").Copy()
$ws2.Range("AH1").PasteSpecial(-4122)
$ws2.Range("AH1").Value = "28. 9. 2021"

# Respondent counts for rows 2..26 for the new wave
$ws2Values = @{
    2  = 1855
    3  = 456
    4  = 678
    5  = 721
    6  = 317
    7  = 325
    8  = 1213
    9  = 897
    10 = 958
    11 = 965
    12 = 429
    13 = 216
    14 = 245
    15 = 39
    16 = 145
    17 = 101
    18 = 21
    19 = 250
    20 = 497
    21 = 237
    22 = 343
    23 = 317
    24 = 219
    25 = 336
    26 = 403
}

for ($row = 2; $row -le 26; $row++) {
    $ws2.Range("AH$row").Value = $ws2Values[$row]
}

# Row 27 is the footnote row; AG27 is an (empty) blank cell that also needs to
# extend into AH27 so the used range covers the new column.
$ws2.Range("AG27").Copy()
$ws2.Range("AH27").PasteSpecial(-4122)

# Update the footnote text in row 27 with the new update date
$ws2.Range("A27").Value = "Život během pandemie, Obavy z epidemie, velikost dotázaného souboru celkově a ve skupinách, aktualizace 6. 10. 2021"
